# Meta VoltageTap A.xlsx — VT interconnect cable specs and w5500 test board
# Applies the content/formatting changes to the "notes" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("notes")

# --- New column D width (connector-picker / links column) ---
$ws.Columns.Item(4).ColumnWidth = 44.3

# --- Row 12: mouser connector picker link ---
$ws.Range("D12").Value = "cool visual connector picker on mouser"
$ws.Range("E12").Value = "https://www.mouser.com/c/i/?number%20of%20positions=8%20Position&packaging=Tube&pitch=2.54%20mm&srsltid=AfmBOoo3AoR97TdMCuXKaB6ufQzjQTeWPaMW7V2CnTmUIDPGfkNKPkwc"

# --- Rows 14-16: merged legal quote block (EV rule) ---
$ws.Range("A14:B16").Merge()
$ws.Range("A14").Value = "EV.5.2.5 Each wire used in a Tractive Battery Container, whether it is part of the GLV or Tractive System, must be rated to the maximum Tractive System voltage"
$ws.Range("A14:B16").HorizontalAlignment = -4108
$ws.Range("A14:B16").VerticalAlignment = -4108
$ws.Range("A14:B16").WrapText = $true
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 30

# --- Row 17: rationale notes ---
$ws.Range("A17").Value = "the intention is that wires don't short"
$ws.Range("B17").Value = "so we can have under specsed wires and connectors, as long as the cable length is 600v rated,, so heat shrink it"
$ws.Rows.Item(17).RowHeight = 45

# --- Rows 20-25: connector / cable part table ---
$ws.Range("A20").Value = "ribbon 8p 2x4 2.54mmp non latched, generic"
$ws.Range("B20").Value = "https://www.mouser.com/ProductDetail/Adam-Tech/FCS-08-SG?qs=xBpwZ0JX2zLXsIHjtia0ew%3D%3D&srsltid=AfmBOoqg4eYgBt3YLFYUVTLoOpNtyC3yUxVKxfLNsrAAmEa31oEmFvKO"
$ws.Rows.Item(20).RowHeight = 75

$ws.Range("A21").Value = "10p 2.54mmp header latched wurth"
$ws.Range("B21").Value = "https://www.mouser.com/ProductDetail/Wurth-Elektronik/61201022121?qs=ZtY9WdtwX54w6hXXLKx1qQ%3D%3D&srsltid=AfmBOoohHCKYZ_5nThgnYgzjyJ7LWaEFKzoos-ZMx7nONpMAa6rFnMsR"
$ws.Rows.Item(21).RowHeight = 60

$ws.Range("A22").Value = "50ft 8pin ribbon cable"
$ws.Range("B22").Value = "https://www.amazon.com/Pc-Accessories-Length-Conductors-Connectors/dp/B00E9P0F34"
$ws.Rows.Item(22).RowHeight = 30

$ws.Range("A23").Value = "8p 2 wrap arround connector"
$ws.Range("B23").Value = "https://www.mouser.com/ProductDetail/Wurth-Elektronik/61200823021?qs=PhR8RmCirEZvQm5v3EiYrA%3D%3D&srsltid=AfmBOooMfGfO819iuGYG45dUWZESwH_55cLvg2yZ3m_V5O8BQQO5v-he"
$ws.Rows.Item(23).RowHeight = 60

$ws.Range("A24").Value = "8p header right angle. Generic. Non latched"
$ws.Range("B24").Value = "https://www.mouser.com/ProductDetail/Wurth-Elektronik/61200821721?qs=PhR8RmCirEb56BXUsQNR%2FQ%3D%3D&srsltid=AfmBOooky3YWgEiw2M0jpH08etGmhJlXVmyZvnTkAkIslv9wCu9MU5VM"
$ws.Rows.Item(24).RowHeight = 60

$ws.Range("A25").Value = "8p ribbon 1.25mmp 1m"
$ws.Range("B25").Value = "https://www.mouser.com/ProductDetail/Wurth-Elektronik/63910815521CAB?qs=rrS6PyfT74fy4EdvcoRmTQ%3D%3D&srsltid=AfmBOool47vkGpVeAHucxC0no_zWHdNkYm7E4z3ruE5nRaMYKgmp8BxC"

# --- View state: scroll so the new rows are visible, select B26 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("B26").Select()
